$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- Header row (row 1): new columns H (date), I (legislator_name), J (legislator_id) ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows 2-12 ---
# Force column H to text so the "2011-12-31" literal isn't auto-converted to a date serial.
$ws.Range("H2:H12").NumberFormat = "@"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = "2011-12-31"
    $ws.Cells.Item($r, 9).Value = "羅明才"
    $ws.Cells.Item($r, 10).Value = 879
}

# --- Replicate formatting from the existing columns onto the new ones ---
# Header style (bold, bordered) from an existing header cell (G1).
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# Data-row style from an existing data cell (G2) onto the new data cells.
$ws.Range("G2").Copy()
$ws.Range("H2:J12").PasteSpecial(-4122)

$excel.CutCopyMode = $false
